$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The F column holds plain "YYYY-MM-DD" text (not real Excel dates) in the
# source data. Excel auto-converts such literals to date serials on
# assignment, so pre-format those cells as Text, write the values, then
# restore the Normal style afterwards so no stray number format survives.
$ws.Range("F2:F9").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "sbkuzh"
$ws.Range("B2").Value = "task_2050-01-01_UZH_LARGE_READY"
$ws.Range("C2").Value = "2024-07-04 22:03:05"
# D2 already empty - left untouched
# E2 already empty - left untouched
$ws.Range("F2").Value = "2050-01-01"
$ws.Range("G2").Value = "LARGE"
$ws.Range("H2").Value = "READY"
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = "sbkrzs"
$ws.Range("B3").Value = "task_2034-01-01_RZS_LARGE_ERROR"
# C3 already empty - left untouched
# D3 already empty - left untouched
# E3 already empty - left untouched
$ws.Range("F3").Value = "2034-01-01"
$ws.Range("G3").Value = "LARGE"
$ws.Range("H3").Value = "ERROR"
$ws.Range("I3").Value = "Missing file task_2034-01-01_RZS_LARGE.xlsx"

# Row 4
$ws.Range("A4").Value = "sbkrzs"
$ws.Range("B4").Value = "task_2033-01-01_RZS_LARGE_ERROR"
$ws.Range("C4").Value = ""
$ws.Range("C4").Style = "Normal"
# D4 already empty - left untouched
$ws.Range("E4").Value = ""
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "2033-01-01"
$ws.Range("G4").Value = "LARGE"
$ws.Range("H4").Value = "ERROR"
$ws.Range("I4").Value = "Missing file task_2033-01-01_RZS_LARGE.xlsx"

# Row 5
$ws.Range("A5").Value = "sbkzbs"
$ws.Range("B5").Value = "task_2033-01-01_ZBS_LARGE_ERROR"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "2033-01-01"
$ws.Range("G5").Value = "LARGE"
$ws.Range("H5").Value = "ERROR"
$ws.Range("I5").Value = "A large task is already scheduled for this date"

# Row 6
$ws.Range("A6").Value = "sbkrzs"
$ws.Range("B6").Value = "task_2041-01-01_UBS_LARGE_ERROR"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "2041-01-01"
$ws.Range("G6").Value = "LARGE"
$ws.Range("H6").Value = "ERROR"
$ws.Range("I6").Value = "Missing file task_2041-01-01_UBS_LARGE.xlsx"

# Row 7
$ws.Range("A7").Value = "sbkrzs"
$ws.Range("B7").Value = "task_2032-01-01_RZS_LARGE_ERROR"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "2032-01-01"
$ws.Range("G7").Value = "LARGE"
$ws.Range("H7").Value = "ERROR"
$ws.Range("I7").Value = "Missing file task_2032-01-01_RZS_LARGE.xlsx"

# Row 8
$ws.Range("A8").Value = "sbkzbz"
$ws.Range("B8").Value = "task_2024-07-06_ZBZ_SMALL_ERROR"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "2024-07-06"
$ws.Range("G8").Value = "SMALL"
$ws.Range("H8").Value = "ERROR"
$ws.Range("I8").Value = "Missing file task_2024-07-06_ZBZ_SMALL.xlsx"

# Row 9
$ws.Range("A9").Value = "sbkhsg"
$ws.Range("B9").Value = "task_2024-07-04_HSG_SMALL_DONE"
$ws.Range("C9").Value = "2024-07-04 22:04:06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "2024-07-04 22:04:23"
$ws.Range("F9").Value = "2024-07-04"
$ws.Range("G9").Value = "SMALL"
$ws.Range("H9").Value = "DONE"
$ws.Range("I9").Style = "Normal"

# Restore the default look of the date column now that the text values are committed.
$ws.Range("F2:F9").Style = "Normal"
